# "Load page position support"
# Adds a second worksheet ("Plan2") that duplicates the layout/data of
# "Plan1" and appends one extra data row (row 8), then repositions the
# selection on both sheets (Plan1 -> F7, Plan2 -> E7) with Plan1 staying
# the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate Plan1 right after itself -> keeps columns/styles/formulas
# (including the F3:F7 shared formula group) identical to the source sheet.
$ws1.Copy($null, $ws1) | Out-Null

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Plan2"

# New row 8: same shape as row 7 (Id=7, Name 6, Prop=9, same Birth Date),
# but its own Value/Calculated Prop.
$ws2.Range("A8").Value = 7
$ws2.Range("B8").Value = "Name 6"
$ws2.Range("C8").Value = 9
$ws2.Range("D8").Value = 31845

# Pull the Birth Date number format from the cell above so D8 renders as
# a date instead of a bare serial number (column D has no default style).
$ws2.Range("D7").Copy() | Out-Null
$ws2.Range("D8").PasteSpecial(-4122) | Out-Null
$ws2.Range("D8").Value = 31845

$ws2.Range("E8").Value = 7

# Re-establish the shared formula for F3:F7 (the copy already carries it,
# this keeps it explicit/intact) and give row 8 its own Calculated Prop
# formula.
$ws2.Range("F3:F7").Formula = "=E3+10"
$ws2.Range("F8").Formula = "=E8+10"

# Selection: Plan2 sits on E7, Plan1 (still the active/visible tab) on F7.
$ws2.Range("E7").Select() | Out-Null
$ws1.Range("F7").Select() | Out-Null
$ws1.Activate() | Out-Null
